# Updates cryptos list values (prices + 1h volume %) to match the latest
# scrape, and fixes the ordering for two coin pairs (dogwifhat/TheGraph and
# FirstDigitalUSD/LidoDAOToken) whose rows were swapped in this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.439.13"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.607.93"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'579.86"
$ws.Range("D6").Value = "'189.75"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("D7").Value = "3.603.28"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +4.22%  "
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "'55.93"
$ws.Range("E12").Value = "  -4.92%  "
$ws.Range("D13").Value = "'0.0000311"
$ws.Range("E13").Value = "  +6.76%  "
$ws.Range("D14").Value = "'9.68"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "4.187.51"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "'19.78"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "3.603.02"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "70.404.13"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "'12.65"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").Value = "'487.36"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'19.39"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("E24").Value = "  -8.66%  "
$ws.Range("D25").Value = "'96.60"
$ws.Range("E25").Value = "  +5.65%  "
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("E27").Value = "  -5.98%  "
$ws.Range("D28").Value = "'11.04"
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("D29").Value = "'9.41"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").Value = "'32.23"
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("D31").Value = "'7.64"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("D32").Value = "'12.18"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "'65.81"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("D35").Value = "'572.50"
$ws.Range("E35").Value = "  -8.65%  "
$ws.Range("D36").Value = "'38.31"
$ws.Range("E36").Value = "  -6.81%  "
$ws.Range("D37").Value = "0.0₃0811"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.34"
$ws.Range("E39").Value = "  +16.13%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.395"
$ws.Range("E40").Value = "  -4.45%  "
$ws.Range("D41").Value = "'2.96"
$ws.Range("E41").Value = "  +5.43%  "
$ws.Range("D42").Value = "'3.50"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("E43").Value = "  -5.97%  "
$ws.Range("D44").Value = "'3.02"
$ws.Range("E44").Value = "  -3.89%  "
$ws.Range("D45").Value = "'3.53"
$ws.Range("E45").Value = "  +6.82%  "
$ws.Range("D46").Value = "3.225.05"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("E48").Value = "  +6.04%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.23"
$ws.Range("E51").Value = "  -3.54%  "
